# Apply crypto price/volume updates to Sheet1 (columns D and E), rows 2-51.
# Values are plain text in the sheet (t="inlineStr"), so numeric-looking
# D-column values must be forced to Text format before assignment to avoid
# Excel auto-converting them to numbers (matches the original text storage).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.407.37"
$ws.Range("E2").Value = "  +0.66%  "
$ws.Range("D3").Value = "3.083.78"
$ws.Range("E3").Value = "  +4.08%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.75"
$ws.Range("E5").Value = "  +0.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.82"
$ws.Range("E6").Value = "  +4.07%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "3.079.92"
$ws.Range("E8").Value = "  +4.03%  "
$ws.Range("E9").Value = "  +1.29%  "
$ws.Range("E10").Value = "  +0.66%  "
$ws.Range("E11").Value = "  +1.10%  "
$ws.Range("E12").Value = "  +5.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000250"
$ws.Range("E13").Value = "  +1.80%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.70"
$ws.Range("E14").Value = "  +7.12%  "
$ws.Range("D16").Value = "3.591.19"
$ws.Range("E16").Value = "  +3.94%  "
$ws.Range("D17").Value = "66.381.12"
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("E18").Value = "  +3.71%  "
$ws.Range("D19").Value = "3.083.99"
$ws.Range("E19").Value = "  +4.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.11"
$ws.Range("E20").Value = "  +17.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "464.57"
$ws.Range("E21").Value = "  +3.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.712"
$ws.Range("E22").Value = "  +5.46%  "
$ws.Range("E23").Value = "  +3.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.19"
$ws.Range("E24").Value = "  +1.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.83"
$ws.Range("E25").Value = "  +5.34%  "
$ws.Range("E26").Value = "  +2.77%  "
$ws.Range("E27").Value = "  +1.28%  "
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.41"
$ws.Range("E30").Value = "  +1.32%  "
$ws.Range("E31").Value = "  +3.22%  "
$ws.Range("E32").Value = "  +2.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.35"
$ws.Range("E33").Value = "  +4.13%  "
$ws.Range("E34").Value = "  +5.50%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("E36").Value = "  +2.27%  "
$ws.Range("E37").Value = "  +2.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.90"
$ws.Range("E38").Value = "  +12.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.09"
$ws.Range("E39").Value = "  +1.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.313"
$ws.Range("E41").Value = "  +1.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.90"
$ws.Range("E42").Value = "  +2.50%  "
$ws.Range("E43").Value = "  +2.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.66"
$ws.Range("E44").Value = "  +3.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0360"
$ws.Range("E45").Value = "  +1.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "383.37"
$ws.Range("E46").Value = "  -0.63%  "
$ws.Range("D47").Value = "2.775.34"
$ws.Range("E47").Value = "  +2.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.68"
$ws.Range("E48").Value = "  +2.83%  "
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.50"
$ws.Range("E50").Value = "  +6.03%  "
$ws.Range("E51").Value = "  +4.56%  "
